# refactor: delete language entity
#
# The "English" language entity is renamed to "ENG":
#   - the worksheet "language_English" becomes "language_ENG"
#   - the settings value "English" (settings!A2, the "default language name")
#     becomes "ENG"
#   - the workbook's active tab/selection moves onto the "settings" sheet,
#     with A2 selected (where the renamed value now lives)

$wb = $excel.ActiveWorkbook

# Rename the worksheet "language_English" -> "language_ENG"
$langSheet = $wb.Worksheets.Item("language_English")
$langSheet.Name = "language_ENG"

# Update the default language name setting from "English" to "ENG"
$settingsSheet = $wb.Worksheets.Item("settings")
$settingsSheet.Range("A2").Value = "ENG"

# The settings sheet becomes the active tab, with A2 as the active selection
$settingsSheet.Activate()
$settingsSheet.Range("A2").Select()
